{"js": "// Replace the date line and each \"a\u00f7b=\" exercise text with its new value.\n// Every \"before\" string in this document is unique, so a simple\n// search-and-replace keyed on the exact old text is unambiguous.\nconst replacements = [\n  [\"2025-05-28 Wednesday\", \"2025-05-29 Thursday\"],\n  [\"594\u00f76=\", \"237\u00f78=\"],\n  [\"811\u00f79=\", \"233\u00f74=\"],\n  [\"655\u00f76=\", \"978\u00f79=\"],\n  [\"959\u00f79=\", \"455\u00f75=\"],\n  [\"825\u00f74=\", \"181\u00f73=\"],\n  [\"555\u00f78=\", \"128\u00f76=\"],\n  [\"752\u00f77=\", \"342\u00f76=\"],\n  [\"538\u00f75=\", \"355\u00f72=\"],\n  [\"729\u00f75=\", \"279\u00f73=\"],\n  [\"177\u00f77=\", \"934\u00f75=\"],\n  [\"410\u00f75=\", \"422\u00f76=\"],\n  [\"124\u00f74=\", \"914\u00f72=\"],\n  [\"387\u00f75=\", \"203\u00f76=\"],\n  [\"513\u00f75=\", \"695\u00f79=\"],\n  [\"658\u00f75=\", \"304\u00f72=\"],\n  [\"820\u00f73=\", \"498\u00f76=\"],\n  [\"412\u00f76=\", \"713\u00f73=\"],\n  [\"117\u00f76=\", \"447\u00f79=\"],\n  [\"781\u00f76=\", \"718\u00f75=\"],\n  [\"832\u00f79=\", \"228\u00f73=\"],\n  [\"358\u00f74=\", \"360\u00f78=\"],\n  [\"264\u00f73=\", \"239\u00f78=\"],\n  [\"312\u00f79=\", \"996\u00f79=\"],\n  [\"849\u00f76=\", \"912\u00f76=\"],\n  [\"921\u00f72=\", \"432\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"a\u00f7b=\" exercise text with its new value.\n# Every \"before\" string in this document is unique, so a simple\n# Find/Replace keyed on the exact old text is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-28 Wednesday\", \"2025-05-29 Thursday\"),\n    @(\"594\u00f76=\", \"237\u00f78=\"),\n    @(\"811\u00f79=\", \"233\u00f74=\"),\n    @(\"655\u00f76=\", \"978\u00f79=\"),\n    @(\"959\u00f79=\", \"455\u00f75=\"),\n    @(\"825\u00f74=\", \"181\u00f73=\"),\n    @(\"555\u00f78=\", \"128\u00f76=\"),\n    @(\"752\u00f77=\", \"342\u00f76=\"),\n    @(\"538\u00f75=\", \"355\u00f72=\"),\n    @(\"729\u00f75=\", \"279\u00f73=\"),\n    @(\"177\u00f77=\", \"934\u00f75=\"),\n    @(\"410\u00f75=\", \"422\u00f76=\"),\n    @(\"124\u00f74=\", \"914\u00f72=\"),\n    @(\"387\u00f75=\", \"203\u00f76=\"),\n    @(\"513\u00f75=\", \"695\u00f79=\"),\n    @(\"658\u00f75=\", \"304\u00f72=\"),\n    @(\"820\u00f73=\", \"498\u00f76=\"),\n    @(\"412\u00f76=\", \"713\u00f73=\"),\n    @(\"117\u00f76=\", \"447\u00f79=\"),\n    @(\"781\u00f76=\", \"718\u00f75=\"),\n    @(\"832\u00f79=\", \"228\u00f73=\"),\n    @(\"358\u00f74=\", \"360\u00f78=\"),\n    @(\"264\u00f73=\", \"239\u00f78=\"),\n    @(\"312\u00f79=\", \"996\u00f79=\"),\n    @(\"849\u00f76=\", \"912\u00f76=\"),\n    @(\"921\u00f72=\", \"432\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)  # wdReplaceAll\n}\n"}
